# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2:73) holds a "quarter date" used to look up the
# matching forecast column (B:AR) for that row. Every value in the
# column is the 1st-of-month serial date for a quarter, but it needs to
# line up with the *next* month's 15th (the mid-month data-release
# convention used elsewhere in the workbook) instead of the quarter's
# first day. Re-stamp each date as "one month later, on the 15th".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's (1900-system) date epoch: serial 0 == 1899-12-30.
$epoch = (Get-Date -Year 1899 -Month 12 -Day 30).Date

for ($row = 2; $row -le 73; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $serial = $cell.Value2

    # Serial -> calendar date (strip any time-of-day noise).
    $oldDate = ($epoch.AddDays($serial)).Date

    # Advance one month, then pin the day-of-month to the 15th.
    $bumped = $oldDate.AddMonths(1)
    $newDate = (Get-Date -Year $bumped.Year -Month $bumped.Month -Day 15).Date

    # Calendar date -> serial, write back.
    $cell.Value = $newDate.ToOADate()
}

Write-Output "Re-indexed column A (rows 2-73) to next-month-15th dates"
